{"js": "// Update the Std_Error, CI_Lower_95, and CI_Upper_95 columns (columns 2-4,\n// 0-indexed) for the five data rows of the ATT estimates table, leaving the\n// Row and ATT_pp columns untouched.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Mapping of old -> new text, keyed by (rowIndex, colIndex) in the table.\n// Row 0 is the header row; data rows are 1-5. Columns: 0=Row, 1=ATT_pp,\n// 2=Std_Error, 3=CI_Lower_95, 4=CI_Upper_95.\nconst updates = [\n  { row: 1, col: 2, text: \"3.15\" },\n  { row: 1, col: 3, text: \"-10.93\" },\n  { row: 1, col: 4, text: \"1.42\" },\n\n  { row: 2, col: 2, text: \"4.22\" },\n  { row: 2, col: 3, text: \"-14.66\" },\n  { row: 2, col: 4, text: \"1.89\" },\n\n  { row: 3, col: 2, text: \"5.68\" },\n  { row: 3, col: 3, text: \"-20.01\" },\n  { row: 3, col: 4, text: \"2.25\" },\n\n  { row: 4, col: 2, text: \"3.52\" },\n  { row: 4, col: 3, text: \"-5.67\" },\n  { row: 4, col: 4, text: \"8.14\" },\n\n  { row: 5, col: 2, text: \"3.17\" },\n  { row: 5, col: 3, text: \"-12.85\" },\n  { row: 5, col: 4, text: \"-0.40\" },\n];\n\nconst paragraphsByUpdate = updates.map((u) => {\n  const cell = table.getCell(u.row, u.col);\n  const paragraphs = cell.body.paragraphs;\n  paragraphs.load(\"items\");\n  return paragraphs;\n});\nawait context.sync();\n\nupdates.forEach((u, i) => {\n  const paragraph = paragraphsByUpdate[i].items[0];\n  const range = paragraph.getRange();\n  range.insertText(u.text, Word.InsertLocation.replace);\n});\n\nawait context.sync();\n", "ps1": "# Update the Std_Error, CI_Lower_95, and CI_Upper_95 columns for the five\n# data rows of the ATT estimates table (Table 1 in the document), leaving\n# the Row and ATT_pp columns untouched.\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\n# Cell() is 1-indexed (row, column). Row 1 is the header row; data rows are\n# 2-6. Columns: 1=Row, 2=ATT_pp, 3=Std_Error, 4=CI_Lower_95, 5=CI_Upper_95.\n$updates = @(\n    @{ Row = 2; Col = 3; Text = \"3.15\" },\n    @{ Row = 2; Col = 4; Text = \"-10.93\" },\n    @{ Row = 2; Col = 5; Text = \"1.42\" },\n\n    @{ Row = 3; Col = 3; Text = \"4.22\" },\n    @{ Row = 3; Col = 4; Text = \"-14.66\" },\n    @{ Row = 3; Col = 5; Text = \"1.89\" },\n\n    @{ Row = 4; Col = 3; Text = \"5.68\" },\n    @{ Row = 4; Col = 4; Text = \"-20.01\" },\n    @{ Row = 4; Col = 5; Text = \"2.25\" },\n\n    @{ Row = 5; Col = 3; Text = \"3.52\" },\n    @{ Row = 5; Col = 4; Text = \"-5.67\" },\n    @{ Row = 5; Col = 5; Text = \"8.14\" },\n\n    @{ Row = 6; Col = 3; Text = \"3.17\" },\n    @{ Row = 6; Col = 4; Text = \"-12.85\" },\n    @{ Row = 6; Col = 5; Text = \"-0.40\" }\n)\n\nforeach ($u in $updates) {\n    $cell = $tbl.Cell($u.Row, $u.Col)\n    $cell.Range.Text = $u.Text\n}\n"}
